$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number but must remain text (inlineStr),
# matching the original cell formatting. Force text format, set value, then restore
# the "Normal" style so no stray style index is left on the cell.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.92'
$ws.Range("D5").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0588'
$ws.Range("D10").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.78'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.83'
$ws.Range("D18").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.21'
$ws.Range("D23").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.30'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.66'
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.93'
$ws.Range("D27").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0471'
$ws.Range("D30").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.07'
$ws.Range("D34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("D35").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.813'
$ws.Range("D40").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.983'
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.79'
$ws.Range("D45").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.44'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0520'
$ws.Range("D49").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0961'
$ws.Range("D51").Style = "Normal"

# Remaining cells: values that are not ambiguous with numeric auto-detection.
$ws.Range("D2").Value = '26.950.84'
$ws.Range("D3").Value = '1.556.36'
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +2.93%  '
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = '1.778.11'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").Value = '1.556.62'
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").Value = '26.955.28'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("E17").Value = '  +0.50%  '
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("E19").Value = '  +1.72%  '
$ws.Range("E20").Value = '  +1.19%  '
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("E22").Value = '  +0.90%  '
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("E25").Value = '  +1.58%  '
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("E30").Value = '  +2.50%  '
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("E32").Value = '  -0.26%  '
$ws.Range("D33").Value = '1.422.89'
$ws.Range("E33").Value = '  +3.88%  '
$ws.Range("E34").Value = '  +4.01%  '
$ws.Range("E35").Value = '  +3.64%  '
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("E38").Value = '  -0.24%  '
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("E40").Value = '  +0.81%  '
$ws.Range("E42").Value = '  +1.05%  '
$ws.Range("E43").Value = '  +3.10%  '
$ws.Range("E44").Value = '  -0.62%  '
$ws.Range("E45").Value = '  +1.76%  '
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("D47").Value = '1.691.62'
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("E48").Value = '  +2.05%  '
$ws.Range("E49").Value = '  +2.51%  '
$ws.Range("D50").Value = '0.0₇0999'
$ws.Range("E50").Value = '  +2.73%  '
$ws.Range("E51").Value = '  +1.46%  '
